$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the account-statement detail rows (B16:J23) so that each worker's
# rows are grouped together and sorted by period (descending), reflecting
# the removal of the previous EC rows and addition of the new ones.
#
# GLORIA VICTORIA VANEGAS REYES (CC 40987503): periods 2212,2211,2210,2209,2208,2207
# NELSON ENRIQUE PACHECO BOHORQUEZ (CC 1043964778): periods 2208,2207

$tipoDoc = "CC"
$gloriaId = "40987503"
$gloriaName = "GLORIA VICTORIA VANEGAS REYES"
$nelsonId = "1043964778"
$nelsonName = "NELSON ENRIQUE PACHECO BOHORQUEZ"

$rows = @(
    @{ Row = 16; Id = $gloriaId; Name = $gloriaName; Period = "2212" },
    @{ Row = 17; Id = $gloriaId; Name = $gloriaName; Period = "2211" },
    @{ Row = 18; Id = $gloriaId; Name = $gloriaName; Period = "2210" },
    @{ Row = 19; Id = $gloriaId; Name = $gloriaName; Period = "2209" },
    @{ Row = 20; Id = $gloriaId; Name = $gloriaName; Period = "2208" },
    @{ Row = 21; Id = $gloriaId; Name = $gloriaName; Period = "2207" },
    @{ Row = 22; Id = $nelsonId; Name = $nelsonName; Period = "2208" },
    @{ Row = 23; Id = $nelsonId; Name = $nelsonName; Period = "2207" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $tipoDoc
    $ws.Cells.Item($r.Row, 3).Value = $r.Id
    $ws.Cells.Item($r.Row, 4).Value = $r.Name
    $ws.Cells.Item($r.Row, 5).Value = $r.Period
}
